$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price/percentage columns (D and E)
# so Excel does not silently convert them to floating point numbers and lose
# the original formatting (trailing zeros, percent signs, etc).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '328.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.97%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '43.86'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.53%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.475'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-2.20%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07988'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-1.07%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.983'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '4.29%'
$ws.Range('B7').Value = 'BTSEToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '2.574'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-3.61%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9496'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.74%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1120'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-5.33%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1879'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '0.42%'
$ws.Range('B11').Value = 'MCDex'
$ws.Range('C11').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '10.71'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '27.04%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09947'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.19%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.04827'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '13.00%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.1064'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.12%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001279'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.73%'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04076'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-2.86%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.005990'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '2.22%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.365'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-6.07%'
$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.380'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.12%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3465'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.01%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1421'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '3.64%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2549'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '1.00%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.001268'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2.41%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-4.43%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0001201'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '1.58%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0003748'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-6.10%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02567'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '-2.75%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05645'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '3.67%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.007555'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-1.91%'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '0.00%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.007395'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '2.94%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002017'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-0.49%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008613'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-2.49%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00007134'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.18%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000751'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.04%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.003534'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '55.57%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.003724'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '5.43%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002102'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.04%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0002002'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.04%'
